$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 12; this shifts existing rows 12-80 down to 13-81
$ws.Rows.Item(12).Insert()

# Fill in the new row 12 with the latest weekly price record
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(12, 3).Value = "Bíobío"
$ws.Cells.Item(12, 4).Value = 45063
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = 100112013
$ws.Cells.Item(12, 7).Value = "Alcachofa"
$ws.Cells.Item(12, 8).Value = "Española"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 7000
$ws.Cells.Item(12, 12).Value = 7500
$ws.Cells.Item(12, 13).Value = 7250
$ws.Cells.Item(12, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 242
$ws.Cells.Item(12, 17).Value = 30
$ws.Cells.Item(12, 18).Value = "Hortaliza"
